$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.775549333333333
$ws.Range("H2").Value = 20.326648
$ws.Range("I2").Value = 0.5307754563424079
$ws.Range("J2").Value = 0.5307754563424079
$ws.Range("M2").Value = 184.1023456666667
$ws.Range("N2").Value = 552.307037
$ws.Range("O2").Value = 0.9813423747591566
$ws.Range("P2").Value = 0.9813423747591565
$ws.Range("Q2").Value = 1247.394525446886
$ws.Range("R2").Value = 11226.55072902198
$ws.Range("S2").Value = 0.5208724467909336
$ws.Range("T2").Value = 0.5208724467909336

# Row 3
$ws.Range("G3").Value = 6.775549333333333
$ws.Range("H3").Value = 20.326648
$ws.Range("I3").Value = 0.5307754563424079
$ws.Range("J3").Value = 0.5307754563424079
$ws.Range("O3").Value = 0.002303378255889225
$ws.Range("P3").Value = 0.002303378255889224
$ws.Range("Q3").Value = 2.927848119403555
$ws.Range("R3").Value = 26.350633074632
$ws.Range("S3").Value = 0.001222576644898783
$ws.Range("T3").Value = 0.001222576644898783

# Row 4
$ws.Range("G4").Value = 6.775549333333333
$ws.Range("H4").Value = 20.326648
$ws.Range("I4").Value = 0.5307754563424079
$ws.Range("J4").Value = 0.5307754563424079
$ws.Range("M4").Value = 1.367901
$ws.Range("N4").Value = 4.103703
$ws.Range("O4").Value = 0.007291483500193526
$ws.Range("P4").Value = 0.007291483500193526
$ws.Range("Q4").Value = 9.268280708616
$ws.Range("R4").Value = 83.41452637754401
$ws.Range("S4").Value = 0.003870140482228356
$ws.Range("T4").Value = 0.003870140482228356

# Row 5
$ws.Range("G5").Value = 6.775549333333333
$ws.Range("H5").Value = 20.326648
$ws.Range("I5").Value = 0.5307754563424079
$ws.Range("J5").Value = 0.5307754563424079
$ws.Range("M5").Value = 1.700197666666667
$ws.Range("N5").Value = 5.100593
$ws.Range("O5").Value = 0.009062763484760617
$ws.Range("P5").Value = 0.009062763484760615
$ws.Range("Q5").Value = 11.51977316691822
$ws.Range("R5").Value = 103.677958502264
$ws.Range("S5").Value = 0.004810292424347127
$ws.Range("T5").Value = 0.004810292424347126

# Row 6
$ws.Range("I6").Value = 0.3421215311185197
$ws.Range("J6").Value = 0.3421215311185197
$ws.Range("M6").Value = 184.1023456666667
$ws.Range("N6").Value = 552.307037
$ws.Range("O6").Value = 0.9813423747591566
$ws.Range("P6").Value = 0.9813423747591565
$ws.Range("Q6").Value = 804.0321379883871
$ws.Range("R6").Value = 7236.289241895484
$ws.Range("S6").Value = 0.3357383558040868
$ws.Range("T6").Value = 0.3357383558040868

# Row 7
$ws.Range("I7").Value = 0.3421215311185197
$ws.Range("J7").Value = 0.3421215311185197
$ws.Range("O7").Value = 0.002303378255889225
$ws.Range("P7").Value = 0.002303378255889224
$ws.Range("S7").Value = 0.0007880352956499271
$ws.Range("T7").Value = 0.0007880352956499269

# Row 8
$ws.Range("I8").Value = 0.3421215311185197
$ws.Range("J8").Value = 0.3421215311185197
$ws.Range("M8").Value = 1.367901
$ws.Range("N8").Value = 4.103703
$ws.Range("O8").Value = 0.007291483500193526
$ws.Range("P8").Value = 0.007291483500193526
$ws.Range("Q8").Value = 5.974048628244
$ws.Range("R8").Value = 53.766437654196
$ws.Range("S8").Value = 0.002494573499211632
$ws.Range("T8").Value = 0.002494573499211632

# Row 9
$ws.Range("I9").Value = 0.3421215311185197
$ws.Range("J9").Value = 0.3421215311185197
$ws.Range("M9").Value = 1.700197666666667
$ws.Range("N9").Value = 5.100593
$ws.Range("O9").Value = 0.009062763484760617
$ws.Range("P9").Value = 0.009062763484760615
$ws.Range("Q9").Value = 7.425291405075111
$ws.Range("R9").Value = 66.827622645676
$ws.Range("S9").Value = 0.003100566519571314
$ws.Range("T9").Value = 0.003100566519571313

# Row 10
$ws.Range("G10").Value = 1.622518
$ws.Range("H10").Value = 4.867554
$ws.Range("I10").Value = 0.1271030125390725
$ws.Range("J10").Value = 0.1271030125390725
$ws.Range("M10").Value = 184.1023456666667
$ws.Range("N10").Value = 552.307037
$ws.Range("O10").Value = 0.9813423747591566
$ws.Range("P10").Value = 0.9813423747591565
$ws.Range("Q10").Value = 298.7093696863887
$ws.Range("R10").Value = 2688.384327177498
$ws.Range("S10").Value = 0.1247315721641363
$ws.Range("T10").Value = 0.1247315721641363

# Row 11
$ws.Range("G11").Value = 1.622518
$ws.Range("H11").Value = 4.867554
$ws.Range("I11").Value = 0.1271030125390725
$ws.Range("J11").Value = 0.1271030125390725
$ws.Range("O11").Value = 0.002303378255889225
$ws.Range("P11").Value = 0.002303378255889224
$ws.Range("Q11").Value = 0.7011219373206667
$ws.Range("R11").Value = 6.310097435886
$ws.Range("S11").Value = 0.0002927663153405151
$ws.Range("T11").Value = 0.000292766315340515

# Row 12
$ws.Range("G12").Value = 1.622518
$ws.Range("H12").Value = 4.867554
$ws.Range("I12").Value = 0.1271030125390725
$ws.Range("J12").Value = 0.1271030125390725
$ws.Range("M12").Value = 1.367901
$ws.Range("N12").Value = 4.103703
$ws.Range("O12").Value = 0.007291483500193526
$ws.Range("P12").Value = 0.007291483500193526
$ws.Range("Q12").Value = 2.219443994718
$ws.Range("R12").Value = 19.974995952462
$ws.Range("S12").Value = 0.000926769518753538
$ws.Range("T12").Value = 0.000926769518753538

# Row 13
$ws.Range("G13").Value = 1.622518
$ws.Range("H13").Value = 4.867554
$ws.Range("I13").Value = 0.1271030125390725
$ws.Range("J13").Value = 0.1271030125390725
$ws.Range("M13").Value = 1.700197666666667
$ws.Range("N13").Value = 5.100593
$ws.Range("O13").Value = 0.009062763484760617
$ws.Range("P13").Value = 0.009062763484760615
$ws.Range("Q13").Value = 2.758601317724667
$ws.Range("R13").Value = 24.827411859522
$ws.Range("S13").Value = 0.001151904540842177
$ws.Range("T13").Value = 0.001151904540842177
